$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as text so the date string is not
# auto-converted into a date serial number, then restore the default
# "Normal" style so no extra number-format style is left on the cell.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "09/13/2025"
$ws.Range("A12").Style = "Normal"

$ws.Range("B12").Value = 0.1167782903132609
$ws.Range("C12").Value = 0.8832217096867391
